{"js": "// Update the date line and the 25 three-digit-divided-by-one-digit\n// expressions in the practice table. Each table cell is addressed by\n// its (row, column) position and the replacement is scoped to that\n// cell's own range, so identical old/new text values occurring in\n// different cells (e.g. \"651\u00f73=\" is both an old value in one cell and\n// the new value written into another) can never cross-contaminate.\n\n// 1) Date paragraph (first paragraph in the body, above the table).\nconst dateResults = context.document.body.search(\"2024-02-10 Saturday\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"2024-02-11 Sunday\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Table of division expressions.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// (row, column, old text, new text) \u2014 row/column are 0-based, matching\n// Table.getCell(row, column). Only the five \"filled\" rows (0, 4, 8, 12,\n// 16) carry text; the intervening rows are blank spacer rows.\nconst replacements = [\n  [0, 0, \"669\u00f76=\", \"854\u00f76=\"],\n  [0, 1, \"393\u00f77=\", \"226\u00f75=\"],\n  [0, 2, \"267\u00f73=\", \"863\u00f73=\"],\n  [0, 3, \"703\u00f76=\", \"972\u00f79=\"],\n  [0, 4, \"467\u00f75=\", \"162\u00f79=\"],\n\n  [4, 0, \"988\u00f77=\", \"651\u00f73=\"],\n  [4, 1, \"578\u00f79=\", \"501\u00f72=\"],\n  [4, 2, \"945\u00f75=\", \"352\u00f76=\"],\n  [4, 3, \"881\u00f75=\", \"953\u00f74=\"],\n  [4, 4, \"651\u00f73=\", \"745\u00f79=\"],\n\n  [8, 0, \"287\u00f79=\", \"896\u00f79=\"],\n  [8, 1, \"758\u00f72=\", \"995\u00f77=\"],\n  [8, 2, \"186\u00f78=\", \"222\u00f78=\"],\n  [8, 3, \"743\u00f72=\", \"915\u00f74=\"],\n  [8, 4, \"236\u00f72=\", \"148\u00f79=\"],\n\n  [12, 0, \"732\u00f74=\", \"706\u00f76=\"],\n  [12, 1, \"686\u00f79=\", \"190\u00f74=\"],\n  [12, 2, \"960\u00f78=\", \"153\u00f77=\"],\n  [12, 3, \"231\u00f78=\", \"139\u00f72=\"],\n  [12, 4, \"849\u00f79=\", \"724\u00f73=\"],\n\n  [16, 0, \"483\u00f73=\", \"196\u00f74=\"],\n  [16, 1, \"195\u00f79=\", \"791\u00f77=\"],\n  [16, 2, \"769\u00f78=\", \"710\u00f75=\"],\n  [16, 3, \"584\u00f72=\", \"213\u00f79=\"],\n  [16, 4, \"189\u00f77=\", \"271\u00f76=\"],\n];\n\nfor (const [row, col, oldText, newText] of replacements) {\n  const cell = table.getCell(row, col);\n  const found = cell.body.search(oldText, { matchCase: true });\n  found.load(\"items\");\n  await context.sync();\n  if (found.items.length > 0) {\n    found.items[0].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 three-digit-divided-by-one-digit\n# expressions in the practice table. Each table cell is addressed by\n# its (row, column) position via $tbl.Cell(row, column) (1-based) and\n# the Find/Replace is scoped to that cell's own Range, so identical\n# old/new text values occurring in different cells (e.g. \"651\u00f73=\" is\n# both an old value in one cell and the new value written into\n# another) can never cross-contaminate.\n\n$d = $word.ActiveDocument\n\n# 1) Date paragraph (first paragraph in the body, above the table).\n$dateRange = $d.Content\n$dateRange.Find.Text = \"2024-02-10 Saturday\"\n$dateFound = $dateRange.Find.Execute()\nif ($dateFound) {\n    $dateRange.Text = \"2024-02-11 Sunday\"\n}\n\n# 2) Table of division expressions.\n$tbl = $d.Tables.Item(1)\n\n# row, column, old text, new text \u2014 row/column are 1-based, matching\n# $tbl.Cell(row, column). Only the five \"filled\" rows (1, 5, 9, 13, 17)\n# carry text; the intervening rows are blank spacer rows.\n$replacements = @(\n    @(1, 1, \"669\u00f76=\", \"854\u00f76=\"),\n    @(1, 2, \"393\u00f77=\", \"226\u00f75=\"),\n    @(1, 3, \"267\u00f73=\", \"863\u00f73=\"),\n    @(1, 4, \"703\u00f76=\", \"972\u00f79=\"),\n    @(1, 5, \"467\u00f75=\", \"162\u00f79=\"),\n\n    @(5, 1, \"988\u00f77=\", \"651\u00f73=\"),\n    @(5, 2, \"578\u00f79=\", \"501\u00f72=\"),\n    @(5, 3, \"945\u00f75=\", \"352\u00f76=\"),\n    @(5, 4, \"881\u00f75=\", \"953\u00f74=\"),\n    @(5, 5, \"651\u00f73=\", \"745\u00f79=\"),\n\n    @(9, 1, \"287\u00f79=\", \"896\u00f79=\"),\n    @(9, 2, \"758\u00f72=\", \"995\u00f77=\"),\n    @(9, 3, \"186\u00f78=\", \"222\u00f78=\"),\n    @(9, 4, \"743\u00f72=\", \"915\u00f74=\"),\n    @(9, 5, \"236\u00f72=\", \"148\u00f79=\"),\n\n    @(13, 1, \"732\u00f74=\", \"706\u00f76=\"),\n    @(13, 2, \"686\u00f79=\", \"190\u00f74=\"),\n    @(13, 3, \"960\u00f78=\", \"153\u00f77=\"),\n    @(13, 4, \"231\u00f78=\", \"139\u00f72=\"),\n    @(13, 5, \"849\u00f79=\", \"724\u00f73=\"),\n\n    @(17, 1, \"483\u00f73=\", \"196\u00f74=\"),\n    @(17, 2, \"195\u00f79=\", \"791\u00f77=\"),\n    @(17, 3, \"769\u00f78=\", \"710\u00f75=\"),\n    @(17, 4, \"584\u00f72=\", \"213\u00f79=\"),\n    @(17, 5, \"189\u00f77=\", \"271\u00f76=\")\n)\n\nforeach ($entry in $replacements) {\n    $row = $entry[0]\n    $col = $entry[1]\n    $oldText = $entry[2]\n    $newText = $entry[3]\n\n    $cellRange = $tbl.Cell($row, $col).Range\n    $cellRange.Find.Text = $oldText\n    $found = $cellRange.Find.Execute()\n    if ($found) {\n        $cellRange.Text = $newText\n    }\n}\n"}
